$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "general" sheet: insert a new "LP solver" row right after the
#    "NLP solver" row (this pushes the existing rows 5-14 down to
#    6-15), fill in its two cells, and give the new label the same
#    style used by every other label in column A.
# ------------------------------------------------------------------
$general = $wb.Worksheets.Item("general")

$general.Rows.Item(5).Insert()

$general.Cells.Item(4, 1).Copy()
$general.Cells.Item(5, 1).PasteSpecial(-4122)

$general.Cells.Item(5, 1).Value = "LP solver (linprog or gurobi)"
$general.Cells.Item(5, 2).Value = "gurobi"

# Widen column A so the new, longer label fits.
$general.Columns.Item(1).ColumnWidth = 48

# ------------------------------------------------------------------
# 2. Make "general" the active / selected sheet (it was previously
#    "measRates"), and move its selection to B6.
# ------------------------------------------------------------------
$general.Activate()
$general.Range("B6").Select()
